# Add four new lead rows (rows 3-6) to the leads export sheet, matching
# the existing row 2 layout/column order (A..V).

$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$rowsData = @(
    @("Prasath", 9876543210, "Inlingua", "Language Center", "Tambaram", "Chennai", "Recurring", "CRM", 100000, $false, 45399, 45423, "Vignesh", "Medium", "prsath@gmail.com", "Call Back", "Nothing", "Nothing", 45401, "admin", "admin", 45399),
    @("Nithya", 9876543210, "Test Company", "IT", "Madipakkam", "Chennai", "Recurring", "Website", 10000, $true, 45399, 45416, "Vignesh", "Medium", "test@abcd.com", "Do Not Disturb", "Nothing", "Nothing", 45416, "admin", "admin", 45400),
    @("Nithya", 9789513442, "SMK", "Politics", "Guindy", "Chennai", "One Time", "SMM", 10000, $true, 45400, 45404, "Chrome", "High", "nithya@gmail.com", "Hold", "abcd", "Nothing", 45396, "admin", "admin", 45400),
    @("Vino", 9456239864, "Abc", "Digital marketing", "Velachery", "Chennai", "One Time", "CRM", 20000, $true, 45400, 45411, "Google", "High", "Vino@gmail.com", "Hold", "cds", "nan", 45405, "admin", "admin", 45400)
)

# Columns that hold dates and need the yyyy-mm-dd number format (same as
# the existing K2/L2/S2/V2 cells).
$dateCols = @("K", "L", "S", "V")

$r = 3
foreach ($row in $rowsData) {
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range($cols[$i] + $r).Value = $row[$i]
    }
    foreach ($dc in $dateCols) {
        $ws.Range($dc + $r).NumberFormat = "yyyy-mm-dd"
    }
    $r++
}
